# Append 5 new subject rows (s46..s50) to the manifest sheet, mirroring the
# existing "meltpatch" feedback rows already present (rows 2-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 47

$data = @(
    @("s46", "s46_e66_321_1_3.jpeg", "meltpatch", "1824", "1755", "77", "38", "88",  "2"),
    @("s47", "s47_e69_321_2_0.jpeg", "meltpatch", "1545", "50",   "77", "38", "129", "2"),
    @("s48", "s48_e64_321_1_1.jpeg", "meltpatch", "1329", "508",  "77", "38", "6",   "2"),
    @("s49", "s49_e70_321_3_1.jpeg", "meltpatch", "1174", "281",  "77", "38", "11",  "2"),
    @("s50", "s50_e67_321_2_2.jpeg", "meltpatch", "1537", "1957", "77", "38", "53",  "2")
)

$endRow = $startRow + $data.Length - 1

# The numeric-looking columns (D:I) must stay stored as text (matching the
# rest of the sheet, which uses inline/text strings rather than numbers), so
# format them as Text before writing the values.
$ws.Range("D" + $startRow + ":I" + $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
